# Corrections to before appendices - M Barnes
#
# 1) The auto "datetimeFigureOut" date field that appears on the Slide
#    Master and on every Custom Layout's Date Placeholder is re-cached
#    from 15/02/2013 to 06/03/2013.
# 2) On the one real slide, the label "Ceramic capacitor plate" is
#    corrected to "Ceramic capacitor" (word "plate" removed), ending up
#    split across two runs ("Ceramic " + "capacitor").

$p = $ppt.ActivePresentation

$oldDate = "15/02/2013"
$newDate = "06/03/2013"

function Update-DatePlaceholders($shapes) {
    $updated = 0
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.Type -eq 14) {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                        $shp.TextFrame.TextRange.Text = $newDate
                        $updated = $updated + 1
                    }
                }
            }
        }
    }
    return $updated
}

# -- Slide Master's own Date Placeholder --
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes | Out-Null

# -- Every Custom Layout's Date Placeholder --
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes | Out-Null
}

# -- Fix the "Ceramic capacitor plate" label on slide 1 --
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Ceramic capacitor plate") {
            # Drop the trailing " plate".
            $tail = $tr.Characters(18, 6)
            $tail.Text = ""

            # Re-set "capacitor" in place so it becomes its own run,
            # matching the two-run split in the target file.
            $tr2 = $shp.TextFrame.TextRange
            $secondRun = $tr2.Characters(9, 9)
            $secondRun.Text = "capacitor"
        }
    }
}
